$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.388.49'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '1.693.90'
$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.33%  '

$ws.Range("D5").Value = '219.25'
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").Value = '0.5501'
$ws.Range("E6").Value = '  +4.32%  '

$ws.Range("D7").Value = '1.010'
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("D8").Value = '0.2740'
$ws.Range("E8").Value = '  +1.24%  '

$ws.Range("D9").Value = '0.06465'
$ws.Range("E9").Value = '  +0.43%  '

$ws.Range("D10").Value = '22.01'
$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("D11").Value = '0.07665'
$ws.Range("E11").Value = '  +2.50%  '

$ws.Range("D12").Value = '1.707.55'
$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("D13").Value = '4.540'
$ws.Range("E13").Value = '  -0.56%  '

$ws.Range("D14").Value = '0.5833'
$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("D15").Value = '0.000008360'
$ws.Range("E15").Value = '  -2.20%  '

$ws.Range("D16").Value = '65.48'
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").Value = '26.432.58'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D20").Value = '10.97'
$ws.Range("E20").Value = '  +0.54%  '

$ws.Range("D21").Value = '191.90'
$ws.Range("E21").Value = '  +0.96%  '

$ws.Range("D22").Value = '6.250'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("D24").Value = '148.84'
$ws.Range("E24").Value = '  +2.99%  '

$ws.Range("D25").Value = '0.1325'
$ws.Range("E25").Value = '  +7.29%  '

$ws.Range("D26").Value = '7.912'
$ws.Range("E26").Value = '  +2.70%  '

$ws.Range("E27").Value = '  -0.99%  '

$ws.Range("D28").Value = '0.06278'
$ws.Range("E28").Value = '  -6.11%  '

$ws.Range("E29").Value = '  +1.73%  '

$ws.Range("D30").Value = '1.333'
$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").Value = '3.599'
$ws.Range("E31").Value = '  +0.25%  '

$ws.Range("D32").Value = '3.607'

$ws.Range("E33").Value = '  +0.66%  '

$ws.Range("E34").Value = '  +1.06%  '

$ws.Range("D35").Value = '0.6148'
$ws.Range("E35").Value = '  -1.39%  '

$ws.Range("D36").Value = '2.413'
$ws.Range("E36").Value = '  +0.76%  '

$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01643'
$ws.Range("E38").Value = '  +1.07%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '6.179'
$ws.Range("E39").Value = '  -3.16%  '

$ws.Range("D40").Value = '1.116.59'
$ws.Range("E40").Value = '  +0.43%  '

$ws.Range("D41").Value = '0.8883'
$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("D43").Value = '101.86'
$ws.Range("E43").Value = '  +0.95%  '

$ws.Range("D44").Value = '1.845.21'
$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000109'
$ws.Range("E45").Value = '  -4.01%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '57.48'
$ws.Range("E46").Value = '  +0.76%  '

$ws.Range("D47").Value = '1.016'
$ws.Range("E47").Value = '  +0.59%  '

$ws.Range("D48").Value = '8.180'
$ws.Range("E48").Value = '  -0.24%  '

$ws.Range("D49").Value = '0.05288'
$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").Value = '6.109'
$ws.Range("E50").Value = '  +0.65%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.4303'
$ws.Range("E51").Value = '  +0.02%  '
